$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "최종점수" (K column) values for rows 2-6
$ws.Range("K2").Value = 59.6
$ws.Range("K3").Value = 58.4
$ws.Range("K4").Value = 50.4
$ws.Range("K5").Value = 49.4
$ws.Range("K6").Value = 46.4

# Update "MACRO_SCORE" (N column) values for rows 2-6
$ws.Range("N2").Value = 54.82400714602223
$ws.Range("N3").Value = 54.82400714602223
$ws.Range("N4").Value = 54.82400714602223
$ws.Range("N5").Value = 54.82400714602223
$ws.Range("N6").Value = 54.82400714602223
